# LOM3231.xlsx edit
#
# The sheet is a label/value form: column A holds a field label (e.g.
# "Programa resumido:") and columns B/C hold the corresponding value
# (duplicated in both columns). The edit removes the standalone
# "519033 - Carlos Yujiro Shigue" row (old row 13, which had no label in
# column A) and shifts every row below it up by one. That shift alone
# reuses several old values one slot earlier than where they belong, so
# the cells are re-filled with their correct values (several of them by
# reusing text that already exists elsewhere on the sheet) before the
# now-empty separator row is deleted.
#
# We do the content copies first (while the original row numbers are
# still valid), working from the bottom of the affected block upward so
# that we never overwrite a cell before its own old value has been
# copied onward. Only then do we delete the old row 13, which shifts
# everything below it up and shrinks the sheet from A1:C25 to A1:C24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 ("Bibliografia:") takes what used to be row 21's value
# ("Aplicação de uma prova..." - the old "Norma de recuperação:" text).
$ws.Range("B21:C21").Copy($ws.Range("B22"))

# Row 21 ("Norma de recuperação:") takes what used to be row 20's value
# ("Média aritmética..." - the old "Critério:" text).
$ws.Range("B20:C20").Copy($ws.Range("B21"))

# Row 20 ("Critério:") takes what used to be row 19's value
# ("Experimentos desenvolvidos..." - the old "Método:" text).
$ws.Range("B19:C19").Copy($ws.Range("B20"))

# Row 19 ("Método:") takes the "519033 - Carlos Yujiro Shigue" value
# (the docente responsável) that used to live in the standalone row 13.
$ws.Range("B13:C13").Copy($ws.Range("B19"))

# Row 10 ("Objetivos:") also now shows the "519033 - Carlos Yujiro
# Shigue" value in place of the old objectives paragraph.
$ws.Range("B13:C13").Copy($ws.Range("B10"))

# Row 16 ("Programa:") now shows the activation date value, copied from
# row 8 ("Ativação:").
$ws.Range("B8:C8").Copy($ws.Range("B16"))

# Row 14 ("Programa resumido:") gets the brand-new value "Semestral".
$ws.Range("B14").Value = "Semestral"
$ws.Range("C14").Value = "Semestral"

# Finally, remove the now-superseded standalone row 13
# ("519033 - Carlos Yujiro Shigue" with a blank label), shifting
# everything below it up by one row.
$ws.Rows.Item(13).Delete()

Write-Output "done"
